# Natmi following Dr Hou advice
# Update the NATMI LR-pair results (Efnb2-Rhbdl2) for rows 2-7 with the
# recomputed ligand/receptor-expressing cell counts (1 -> 3) and the
# resulting recalculated expression / specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = [ordered]@{
    2 = @{
        E = 3
        G = 31.910369
        H = 95.73110699999999
        I = 0.579978174461379
        J = 0.579978174461379
        K = 3
        M = 1.021610666666667
        N = 3.064832
        Q = 32.59997334766933
        R = 293.399760129024
        S = 0.579978174461379
        T = 0.579978174461379
    }
    3 = @{
        E = 3
        G = 11.420477
        H = 34.261431
        I = 0.2075697527013294
        J = 0.2075697527013294
        K = 3
        M = 1.021610666666667
        N = 3.064832
        Q = 11.66728112162133
        R = 105.005530094592
        S = 0.2075697527013294
        T = 0.2075697527013294
    }
    4 = @{
        E = 3
        G = 0.6836493333333333
        H = 2.050948
        I = 0.01242548126969028
        J = 0.01242548126969028
        K = 3
        M = 1.021610666666667
        N = 3.064832
        Q = 0.6984234511928888
        R = 6.285811060736
        S = 0.01242548126969028
        T = 0.01242548126969028
    }
    5 = @{
        E = 3
        G = 0.4491346666666667
        H = 1.347404
        I = 0.00816312415756312
        J = 0.00816312415756312
        K = 3
        M = 1.021610666666667
        N = 3.064832
        Q = 0.4588407662364445
        R = 4.129566896128
        S = 0.00816312415756312
        T = 0.00816312415756312
    }
    6 = @{
        E = 3
        G = 0.6568320000000001
        H = 1.970496
        I = 0.01193807017047708
        J = 0.01193807017047708
        K = 3
        M = 1.021610666666667
        N = 3.064832
        Q = 0.6710265774080001
        R = 6.039239196672001
        S = 0.01193807017047708
        T = 0.01193807017047708
    }
    7 = @{
        E = 3
        G = 9.899486
        H = 29.698458
        I = 0.1799253972395612
        J = 0.1799253972395612
        K = 3
        M = 1.021610666666667
        N = 3.064832
        Q = 10.11342049211733
        R = 91.02078442905599
        S = 0.1799253972395612
        T = 0.1799253972395612
    }
}

foreach ($rowNum in $rows.Keys) {
    $colValues = $rows[$rowNum]
    foreach ($col in $colValues.Keys) {
        $ws.Range("$col$rowNum").Value = $colValues[$col]
    }
}
